$d = $word.ActiveDocument

# The table "2.3. Table comentario" is the 3rd table in the document body.
$table = $d.Tables.Item(3)

# Find the row whose first cell contains "usuario_id" so we insert the new
# "estrelas" row right above it (rows: id, texto, usuario_id, jogo_id, hora).
$targetIndex = 0
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    $cellText = $table.Cell($i, 1).Range.Text
    if ($cellText -like "usuario_id*") {
        $targetIndex = $i
        break
    }
}

$targetRow = $table.Rows.Item($targetIndex)
$newRow = $table.Rows.Add($targetRow)
$newIndex = $newRow.Index

$table.Cell($newIndex, 1).Range.Text = "estrelas"
$table.Cell($newIndex, 2).Range.Text = "int"
$table.Cell($newIndex, 3).Range.Text = ""
$table.Cell($newIndex, 4).Range.Text = ""
